$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.049.04"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.61%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.861.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.30%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.15%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.23%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5116"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.76%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3742"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.01%  "

$ws.Range("E9").Value = "  -1.89%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8920"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.96%  "

$ws.Range("E11").Value = "  -2.58%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07542"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.23%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.854.47"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.77%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.299"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.94%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "88.83"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.27%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.0000"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.26%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008369"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.90%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.09"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.06%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.0000"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.19%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.095.62"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.53%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.064"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.83%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.094.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.81%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.54"
$ws.Range("D23").Style = "Normal"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.473"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.64%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.26"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.38%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.840"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.37%  "

$ws.Range("E27").Value = "  -1.90%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.089"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.55%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "112.91"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.82%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.686"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.81%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.656"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.84%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09037"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.12%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05115"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.34%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.050"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.10%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.157"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.94%  "

$ws.Range("E36").Value = "  -6.86%  "

$ws.Range("E37").Value = "  -1.21%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.507"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.78%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.054"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.14%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.071"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.78%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5335"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.15%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.593"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.46%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "115.39"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.77%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.322"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.68%  "

$ws.Range("E45").Value = "  -2.71%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.27%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4625"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.24%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.995"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.31%  "

$ws.Range("E49").Value = "  -3.99%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.77"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.33%  "

$ws.Range("E51").Value = "  -4.38%  "
